$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61: move "Scalpel Accuracy:" label from C61 to E61,
# clear the old numeric value in D61, and put the corrected
# accuracy value in F61.
$ws.Range("C61").Value = $null
$ws.Range("D61").Value = $null
$ws.Range("E61").Value = "Scalpel Accuracy:"
$ws.Range("F61").Value = 68.97

# Row 62: fix label wording
$ws.Range("E62").Value = "Accuracy vs PyType"
